# Fruta / hortaliza, semanal
# The weekly data rows (2-24) get re-shuffled across rows for columns
# D (Fecha), K (Variedad), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen), S (Precio $/Kg), T (Kg / unidad).
# Columns A,B,C,E,F,G,H,I,J are identical on every row and are left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values currently sitting in
# source row should end up in destination row).
$map = @{
    2  = 4
    3  = 20
    4  = 16
    5  = 15
    6  = 14
    7  = 7
    8  = 11
    9  = 17
    10 = 6
    11 = 8
    12 = 3
    13 = 10
    14 = 12
    15 = 5
    16 = 2
    17 = 13
    18 = 24
    19 = 18
    20 = 19
    21 = 9
    22 = 21
    23 = 23
    24 = 22
}

$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # D,K,L,M,N,O,P,Q,R,S,T

# Snapshot every relevant cell value before writing anything, since the
# mapping moves data between rows and we must not read already-overwritten
# cells.
$snapshot = @{}
for ($r = 2; $r -le 24; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

for ($r = 2; $r -le 24; $r++) {
    $src = $map[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}
